$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking "cryptos" table refresh (GitHub Actions bot, 2024-07-27).
# Every data row keeps its original column layout/styling; only the
# scraped Coin/Link/Price/Volume(1h) text is rewritten cell-by-cell,
# including two header blocks whose rank flipped with a neighbour
# (Dai<->Litecoin, Bittensor/InjectiveProtocol/Maker/Hedera reshuffle,
# Stellar<->ONDO).
#
# A handful of refreshed Price values are plain decimals (e.g. "1.00",
# "0.999") that Excel would otherwise auto-coerce to numbers, losing the
# trailing zeros / changing the stored type away from text. Prefixing
# with a quote forces those specific cells to stay literal text, exactly
# like the rest of the column (the source data is text, never a number).

$ws.Range('D2').Value = '68.131.53'
$ws.Range('E2').Value = '  +1.65%  '
$ws.Range('D3').Value = '3.281.43'
$ws.Range('E3').Value = '  +0.83%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''588.09'
$ws.Range('E5').Value = '  +1.61%  '
$ws.Range('D6').Value = '''186.67'
$ws.Range('E6').Value = '  +5.22%  '
$ws.Range('E8').Value = '  +0.31%  '
$ws.Range('D9').Value = '''0.135'
$ws.Range('E9').Value = '  +4.17%  '
$ws.Range('D10').Value = '''6.73'
$ws.Range('E10').Value = '  -0.19%  '
$ws.Range('D11').Value = '''0.418'
$ws.Range('E11').Value = '  +1.18%  '
$ws.Range('D12').Value = '3.847.53'
$ws.Range('E12').Value = '  +0.66%  '
$ws.Range('E13').Value = '  +0.57%  '
$ws.Range('D14').Value = '''28.71'
$ws.Range('E14').Value = '  +2.30%  '
$ws.Range('D15').Value = '68.174.79'
$ws.Range('E15').Value = '  +1.74%  '
$ws.Range('E16').Value = '  +2.99%  '
$ws.Range('D17').Value = '3.279.05'
$ws.Range('E17').Value = '  +0.81%  '
$ws.Range('E18').Value = '  +0.55%  '
$ws.Range('D19').Value = '''13.69'
$ws.Range('E19').Value = '  +2.10%  '
$ws.Range('D20').Value = '''382.92'
$ws.Range('E20').Value = '  +2.68%  '
$ws.Range('E21').Value = '  +1.41%  '
$ws.Range('B22').Value = 'Litecoin'
$ws.Range('C22').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D22').Value = '''71.56'
$ws.Range('E22').Value = '  +0.55%  '
$ws.Range('B23').Value = 'Dai'
$ws.Range('C23').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D23').Value = '''1.00'
$ws.Range('E23').Value = '  -0.18%  '
$ws.Range('E24').Value = '  +1.05%  '
$ws.Range('E25').Value = '  +2.29%  '
$ws.Range('D26').Value = '''9.91'
$ws.Range('E26').Value = '  +0.85%  '
$ws.Range('D27').Value = '''0.186'
$ws.Range('E27').Value = '  +3.99%  '
$ws.Range('D28').Value = '''0.999'
$ws.Range('E28').Value = '  +0.00%  '
$ws.Range('D29').Value = '''5.86'
$ws.Range('E29').Value = '  +4.69%  '
$ws.Range('D30').Value = '''2.01'
$ws.Range('E30').Value = '  +1.29%  '
$ws.Range('D31').Value = '''7.28'
$ws.Range('E31').Value = '  +6.73%  '
$ws.Range('D32').Value = '''22.94'
$ws.Range('E32').Value = '  +1.48%  '
$ws.Range('E33').Value = '  +2.15%  '
$ws.Range('E34').Value = '  -0.01%  '
$ws.Range('E35').Value = '  +3.20%  '
$ws.Range('D36').Value = '''162.28'
$ws.Range('E36').Value = '  -2.97%  '
$ws.Range('E37').Value = '  +0.23%  '
$ws.Range('D38').Value = '''0.841'
$ws.Range('E38').Value = '  -1.42%  '
$ws.Range('D39').Value = '''6.82'
$ws.Range('E39').Value = '  +5.77%  '
$ws.Range('E40').Value = '  -1.47%  '
$ws.Range('E41').Value = '  +5.65%  '
$ws.Range('D42').Value = '''2.63'
$ws.Range('E42').Value = '  +1.71%  '
$ws.Range('D43').Value = '''41.50'
$ws.Range('E43').Value = '  +2.46%  '
$ws.Range('B44').Value = 'InjectiveProtocol'
$ws.Range('C44').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D44').Value = '''25.58'
$ws.Range('E44').Value = '  +1.22%  '
$ws.Range('B45').Value = 'Maker'
$ws.Range('C45').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D45').Value = '2.658.67'
$ws.Range('E45').Value = '  -3.68%  '
$ws.Range('B46').Value = 'Hedera'
$ws.Range('C46').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D46').Value = '''0.0689'
$ws.Range('E46').Value = '  +1.89%  '
$ws.Range('B47').Value = 'Bittensor'
$ws.Range('C47').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D47').Value = '''345.97'
$ws.Range('E47').Value = '  -1.05%  '
$ws.Range('D48').Value = '''0.0286'
$ws.Range('E48').Value = '  +1.95%  '
$ws.Range('E49').Value = '  +5.22%  '
$ws.Range('B50').Value = 'ONDO'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range('D50').Value = '''1.00'
$ws.Range('E50').Value = '  +2.13%  '
$ws.Range('B51').Value = 'Stellar'
$ws.Range('C51').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D51').Value = '''0.103'
$ws.Range('E51').Value = '  +0.30%  '
